# Battle Log.xlsx revision: "Revised combat,py to use command class"
#
# 1) Group MANA - Enemy!L5 COMMAND changes from "ColtGun" to "Colt".
# 2) A brand-new "Static DMG" sheet is appended to the end of the workbook,
#    mirroring the layout of the other combat sheets (ROBO/Colt vs Goblin).

$wb = $excel.ActiveWorkbook

# --- 1. Rename the ROBO's command on "Group MANA - Enemy" ------------------
$enemyWs = $wb.Worksheets.Item("Group MANA - Enemy")
$enemyWs.Range("L5").Value = "Colt"

# --- 2. Append the new "Static DMG" worksheet -------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "Static DMG"

# Header row (shared across every combat-log sheet in this workbook)
$ws.Cells.Item(1, 1).Value = "Index"
$ws.Cells.Item(1, 2).Value = "NAME"
$ws.Cells.Item(1, 3).Value = "ROLE"
$ws.Cells.Item(1, 4).Value = "LIVES"
$ws.Cells.Item(1, 5).Value = "Position"
$ws.Cells.Item(1, 6).Value = "Initiative"
$ws.Cells.Item(1, 7).Value = "CURRENT HP"
$ws.Cells.Item(1, 8).Value = "CURRENT STR"
$ws.Cells.Item(1, 9).Value = "CURRENT AGL"
$ws.Cells.Item(1, 10).Value = "CURRENT MANA"
$ws.Cells.Item(1, 11).Value = "CURRENT DEF"
$ws.Cells.Item(1, 12).Value = "COMMAND"
$ws.Cells.Item(1, 13).Value = "TARGET"
$ws.Cells.Item(1, 14).Value = "Stoned"
$ws.Cells.Item(1, 15).Value = "Cursed"
$ws.Cells.Item(1, 16).Value = "Blinded"
$ws.Cells.Item(1, 17).Value = "Stunned"
$ws.Cells.Item(1, 18).Value = "Paralyzed"
$ws.Cells.Item(1, 19).Value = "Poisoned"
$ws.Cells.Item(1, 20).Value = "Confused"
$ws.Cells.Item(1, 21).Value = "ACTIONS TAKEN"

# Row 2 - ROBO (Player), using the "Colt" command against the Goblin
$ws.Cells.Item(2, 2).Value = "ROBO"
$ws.Cells.Item(2, 1).Formula = "=B2"
$ws.Cells.Item(2, 3).Value = "Player"
$ws.Cells.Item(2, 4).Value = 1
$ws.Cells.Item(2, 5).Value = 1
$ws.Cells.Item(2, 12).Value = "Colt"
$ws.Cells.Item(2, 13).Value = "Goblin"

# Row 3 - Goblin (Enemy)
$ws.Cells.Item(3, 2).Value = "Goblin"
$ws.Cells.Item(3, 1).Formula = "=B3"
$ws.Cells.Item(3, 3).Value = "Enemy"
$ws.Cells.Item(3, 4).Value = 1

# Column widths matching the other sheets that share this layout
$ws.Range("D1:E1").ColumnWidth = 11.7109375
$ws.Range("G1").ColumnWidth = 12.85546875
$ws.Range("H1").ColumnWidth = 14
$ws.Range("I1").ColumnWidth = 14.140625
$ws.Range("J1").ColumnWidth = 15.85546875
$ws.Range("K1").ColumnWidth = 14.140625
$ws.Range("L1").ColumnWidth = 12.42578125
$ws.Range("U1").ColumnWidth = 15.85546875

# Put the cursor/view where the author left it
$ws.Application.ActiveWindow.TopLeftCell = $ws.Range("B1")
$ws.Range("L3").Select()
